# v2.6 Added decoupled suspension, four-wheel steering, scripts to generate GGV diagram
#
# Adds a new worksheet "Semi_Truck_Scalable" (a scaled copy of
# "Truck_Amandla_3Axle") to the workbook, placed after
# "Truck_Amandla_3Axle", and updates the active-sheet/selection state
# that Excel persists per sheet.

$wb = $excel.ActiveWorkbook

$wsBus   = $wb.Worksheets.Item("Bus_Makhulu_3Axle")
$wsTruck = $wb.Worksheets.Item("Truck_Amandla_3Axle")

# New sheet is a copy of Truck_Amandla_3Axle, inserted right after it.
$wsTruck.Copy($null, $wsTruck)
$wsSemi = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsSemi.Name = "Semi_Truck_Scalable"

# Label the new sheet's "Instance" cell with its own name.
$wsSemi.Range("H3").Value = "Semi_Truck_Scalable"

# Truck_Amandla_3Axle keeps its last selection, now at D24, and is no
# longer the active tab.
$wsTruck.Activate()
$wsTruck.Range("D24").Select()

# The new sheet becomes the active tab, with J17 selected.
$wsSemi.Activate()
$wsSemi.Range("J17").Select()
